$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.224.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5296"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.007"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +0.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06383"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.58"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07673"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.698.92"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.513"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5766"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008346"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.257.12"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.007"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.874"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.91%  "
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "190.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.225"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "148.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.797"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("E26").Value = "  -0.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06249"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.372"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.572"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.560"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.683"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.023"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6147"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.422"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.758"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.180"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01622"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.104.19"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.15%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8902"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.011"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.834.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000111"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.48"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.006"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.069"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05277"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4286"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.17%  "
$ws.Range("E51").Value = "  -0.44%  "
